$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: fill in the two new trailing columns (PriceChange / UpDown) ---
$ws.Range("X3").Value = 4.75
$ws.Range("Y3").Value = "Up"

# --- Row 4: new observation appended to the table ---
# Copy formatting from row 3 so date/percent styles match exactly.
$ws.Range("A3").Copy($ws.Range("A4"))
$ws.Range("A4").Value = 42641.891504629632

$ws.Range("B4").Value = 6
$ws.Range("C4").Value = "Buy"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = "Random"
$ws.Range("Q4").Value = 57.519894101767122
$ws.Range("R4").Value = 1.83

$ws.Range("S3").Copy($ws.Range("S4"))
$ws.Range("S4").Value = 0.1364

$ws.Range("T3").Copy($ws.Range("T4"))
$ws.Range("T4").Value = 0.0165

$ws.Range("U4").Value = 6.04
$ws.Range("V4").Value = "N/A"
$ws.Range("W4").Value = 2

# Note: adding the row nudges Excel's own best-fit column widths by a
# sub-pixel amount in the source workbook, but that's a passive side effect
# of recalculating best-fit (the columns stay bestFit/customWidth either
# way) rather than a deliberate width edit, so it's left for Excel/the
# engine to derive rather than hard-coding fragile pixel-grid values here.
